# Add new certificate row (Git & Github Bootcamp completion) to Sheet1, row 9.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the Date (column B) first as text (not an auto-converted date). A leading
# apostrophe forces Excel to store it as text (reusing the existing "quote
# prefix" style) instead of auto-converting the "2023.10.21"-like string into a
# real date (which would otherwise create a brand-new number-format style).
# Then copy formats back from a neighboring cell that uses the original style
# (s="1") so no new style entries are introduced in styles.xml.
$ws.Range("B9").Value2 = "'2023.10.21"
$ws.Range("A8").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the Cert No. (column A) - plain text, not date-like, so it stays text
# automatically with the original style.
$ws.Range("A9").Value2 = "UC-e26431f6-ca7b-4eb2-827b-9f42a8cd8399"

# Reflect the final cell selection/active cell as recorded in the saved file.
$ws.Range("A12").Select()
